$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "el cual quedo alojado en la carpeta principal.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "el cual quedo alojado en la carpeta principal: sitemap.xml.",
    2
)
